# PIB_potencial_rm.xlsx -- "comparison prod fun , hp, cf lac 18"
#
# Summary of the edit (per the OOXML diff):
#   * Sheet1: add a header label "iso2c" in A2 (column A holds the year axis).
#   * Sheet2:
#       - add a header label "iso2c" in A1
#       - drop the "2003-2008" / "2010-2016" summary headers in AB1:AC1
#       - rename three country-code rows: A9 ESV -> SV, A17 RD -> DO, A18 UR -> UY
#       - clear out the AVERAGE() helper formulas that lived in AB2:AC19
#         (the cells stay, just emptied of their formula/value)
#       - remove the little "lac18 / cardm / sa" legend block that lived in
#         AA21:AA23 below the table
#   * Sheet2 becomes the active sheet/tab (previously Sheet3 was active).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet1: new header cell -------------------------------------------------
$ws1.Range("A2").Value = "iso2c"

# --- Sheet2: header row -------------------------------------------------------
$ws2.Range("A1").Value = "iso2c"
$ws2.Range("AB1:AC1").ClearContents()

# --- Sheet2: country code relabelling ---------------------------------------
$ws2.Range("A9").Value = "SV"
$ws2.Range("A17").Value = "DO"
$ws2.Range("A18").Value = "UY"

# --- Sheet2: drop the AVERAGE() helper columns -------------------------------
$ws2.Range("AB2:AC19").ClearContents()

# --- Sheet2: drop the trailing legend block (AA21:AA23) ----------------------
$ws2.Range("A21:AC23").ClearContents()

# --- Selection / active-sheet bookkeeping (matches the saved view state) ----
$ws1.Activate()
$ws1.Range("A3").Select()

$ws3.Activate()
$ws3.Range("B2").Select()

$ws2.Activate()
$ws2.Range("A17").Select()

Write-Output "applied PIB_potencial_rm edits"
